# Update crypto price/volume table to reflect the latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.724.57"
$ws.Range("E2").Value = "  +0.68%  "

$ws.Range("D3").Value = "1.962.58"
$ws.Range("E3").Value = "  +0.99%  "

$ws.Range("E4").Value = "  -0.10%  "

$ws.Range("D5").Value = "'244.86"
$ws.Range("E5").Value = "  +0.41%  "

$ws.Range("D6").Value = "'0.620"
$ws.Range("E6").Value = "  +1.09%  "

$ws.Range("D7").Value = "'59.03"
$ws.Range("E7").Value = "  +1.28%  "

$ws.Range("E9").Value = "  +1.88%  "

$ws.Range("D10").Value = "'0.0810"
$ws.Range("E10").Value = "  -3.38%  "

$ws.Range("E11").Value = "  -0.39%  "

$ws.Range("D12").Value = "'22.43"
$ws.Range("E12").Value = "  +4.32%  "

$ws.Range("D13").Value = "2.249.27"
$ws.Range("E13").Value = "  +0.92%  "

$ws.Range("D14").Value = "'0.826"
$ws.Range("E14").Value = "  -0.03%  "

$ws.Range("D15").Value = "'13.75"
$ws.Range("E15").Value = "  +1.08%  "

$ws.Range("D16").Value = "'5.28"
$ws.Range("E16").Value = "  +0.43%  "

$ws.Range("D17").Value = "1.967.86"
$ws.Range("E17").Value = "  +1.78%  "

$ws.Range("D18").Value = "36.620.17"
$ws.Range("E18").Value = "  +0.58%  "

$ws.Range("D19").Value = "'69.78"
$ws.Range("E19").Value = "  +0.00%  "

$ws.Range("D20").Value = "0.0₃0862"
$ws.Range("E20").Value = "  -1.09%  "

$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").Value = "'5.09"
$ws.Range("E21").Value = "  +1.78%  "

$ws.Range("B22").Value = "BitcoinCash"
$ws.Range("C22").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D22").Value = "'229.24"
$ws.Range("E22").Value = "  -0.23%  "

$ws.Range("D23").Value = "'1.00"
$ws.Range("E23").Value = "  -0.17%  "

$ws.Range("D24").Value = "'2.44"
$ws.Range("E24").Value = "  -0.60%  "

$ws.Range("D25").Value = "'2.36"
$ws.Range("E25").Value = "  +3.03%  "

$ws.Range("D26").Value = "'9.33"
$ws.Range("E26").Value = "  +0.19%  "

$ws.Range("E27").Value = "  +13.63%  "

$ws.Range("D28").Value = "'160.68"
$ws.Range("E28").Value = "  -1.20%  "

$ws.Range("D29").Value = "'19.42"
$ws.Range("E29").Value = "  -0.16%  "

$ws.Range("E30").Value = "  +1.11%  "

$ws.Range("D31").Value = "'1.13"
$ws.Range("E31").Value = "  -1.96%  "

$ws.Range("D32").Value = "'4.73"
$ws.Range("E32").Value = "  +0.80%  "

$ws.Range("D33").Value = "'0.0619"
$ws.Range("E33").Value = "  -2.18%  "

$ws.Range("D34").Value = "'4.27"
$ws.Range("E34").Value = "  -0.17%  "

$ws.Range("E35").Value = "  +0.03%  "

$ws.Range("D36").Value = "'6.08"
$ws.Range("E36").Value = "  -2.50%  "

$ws.Range("D37").Value = "'2.26"
$ws.Range("E37").Value = "  +4.93%  "

$ws.Range("D38").Value = "'3.40"
$ws.Range("E38").Value = "  +12.22%  "

$ws.Range("E39").Value = "  -0.11%  "

$ws.Range("D40").Value = "'0.100"
$ws.Range("E40").Value = "  +3.38%  "

$ws.Range("E41").Value = "  -2.24%  "

$ws.Range("E42").Value = "  +1.32%  "

$ws.Range("D43").Value = "'1.17"
$ws.Range("E43").Value = "  -1.48%  "

$ws.Range("D44").Value = "'16.11"
$ws.Range("E44").Value = "  +0.16%  "

$ws.Range("D45").Value = "1.358.48"
$ws.Range("E45").Value = "  +0.43%  "

$ws.Range("E46").Value = "  +0.35%  "

$ws.Range("D47").Value = "'87.78"
$ws.Range("E47").Value = "  -0.05%  "

$ws.Range("D48").Value = "'7.16"
$ws.Range("E48").Value = "  -0.73%  "

$ws.Range("E49").Value = "  +0.67%  "

$ws.Range("D50").Value = "2.140.81"

$ws.Range("D51").Value = "'43.73"
$ws.Range("E51").Value = "  -4.14%  "
